$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.735.22"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "3.101.29"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "524.56"
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("D6").Value = "141.93"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.101.55"
$ws.Range("E8").Value = "  +1.21%  "
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("E12").Value = "  +3.61%  "
$ws.Range("D13").Value = "3.637.65"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("D14").Value = "0.131"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("D15").Value = "25.78"
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "57.862.79"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").Value = "3.102.21"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").Value = "6.10"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("D22").Value = "341.54"
$ws.Range("E22").Value = "  +3.11%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").Value = "66.90"
$ws.Range("E25").Value = "  +2.33%  "
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "0.0₃0919"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").Value = "6.50"
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("E32").Value = "  +3.62%  "
$ws.Range("D33").Value = "21.03"
$ws.Range("E33").Value = "  +1.89%  "
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("D35").Value = "155.69"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").Value = "4.63"
$ws.Range("E36").Value = "  +2.05%  "
$ws.Range("D37").Value = "6.15"
$ws.Range("E37").Value = "  +2.20%  "
$ws.Range("D38").Value = "27.18"
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("D39").Value = "1.25"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("E40").Value = "  -2.10%  "
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("D42").Value = "3.143.18"
$ws.Range("E42").Value = "  +1.11%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.684"
$ws.Range("E43").Value = "  +4.04%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "1.52"
$ws.Range("E44").Value = "  +9.42%  "
$ws.Range("D45").Value = "36.72"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "2.295.16"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("D49").Value = "0.984"
$ws.Range("E49").Value = "  +4.84%  "
$ws.Range("D50").Value = "20.55"
$ws.Range("E50").Value = "  -1.54%  "
$ws.Range("D51").Value = "6.03"
$ws.Range("E51").Value = "  +1.65%  "
